$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F2").Value = 886
$ws1.Range("F3").Value = 111
$ws1.Range("F5").Value = 2532
$ws1.Range("F6").Value = 681
$ws1.Range("F8").Value = 4208
$ws1.Range("F10").Value = 400
$ws1.Range("F11").Value = 3318
$ws1.Range("F12").Value = 924
$ws1.Range("F15").Value = 252
$ws1.Range("F16").Value = 2239
$ws1.Range("F17").Value = 1232
$ws1.Range("F18").Value = 19
$ws1.Range("F19").Value = 1987
$ws1.Range("F20").Value = 463
$ws1.Range("F22").Value = 43
$ws1.Range("F23").Value = 9265
$ws1.Range("F24").Value = 5856
$ws1.Range("F27").Value = 787
$ws1.Range("F29").Value = 804
$ws1.Range("F30").Value = 3487
$ws1.Range("F32").Value = 955
$ws1.Range("F33").Value = 429
$ws1.Range("F35").Value = 219
$ws1.Range("F36").Value = 190
$ws1.Range("F37").Value = 4741
$ws1.Range("F38").Value = 12
$ws1.Range("F39").Value = 952
$ws2.Range("F11").Value = 85
$ws2.Range("F15").Value = 3499
$ws3.Range("F2").Value = 8551
$ws3.Range("F3").Value = 399
$ws3.Range("F4").Value = 1446
$ws4.Range("F2").Value = 8551
$ws4.Range("F3").Value = 886
$ws4.Range("F4").Value = 399
$ws4.Range("F5").Value = 1446
$ws4.Range("F6").Value = 111
$ws4.Range("F9").Value = 4208
$ws4.Range("F11").Value = 400
$ws4.Range("F12").Value = 3318
$ws4.Range("F13").Value = 924
$ws4.Range("F16").Value = 252
$ws4.Range("F17").Value = 2239
$ws4.Range("F22").Value = 1232
$ws4.Range("F23").Value = 85
$ws4.Range("F24").Value = 19
$ws4.Range("F26").Value = 463
$ws4.Range("F28").Value = 43
$ws4.Range("F29").Value = 9265
$ws4.Range("F30").Value = 3499
$ws4.Range("F35").Value = 804
$ws4.Range("F37").Value = 955
$ws4.Range("F38").Value = 429
$ws4.Range("F40").Value = 219
$ws4.Range("F42").Value = 190
$ws4.Range("F43").Value = 4741
$ws4.Range("F44").Value = 12
$ws4.Range("F45").Value = 952
